$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.383.68"
$ws.Range("E2").Value = "  +0.76%  "
$ws.Range("D3").Value = "'2.582.88"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").Value = "'507.51"
$ws.Range("E5").Value = "  -0.50%  "
$ws.Range("D6").Value = "'153.13"
$ws.Range("E6").Value = "  -3.72%  "
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("D8").Value = "'0.579"
$ws.Range("E8").Value = "  -5.16%  "
$ws.Range("D9").Value = "'2.582.27"
$ws.Range("E9").Value = "  -1.33%  "
$ws.Range("D10").Value = "'6.56"
$ws.Range("E10").Value = "  +6.65%  "
$ws.Range("D11").Value = "'0.104"
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").Value = "'0.348"
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").Value = "'3.022.35"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").Value = "'60.385.67"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").Value = "'21.58"
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("D17").Value = "'0.0000140"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").Value = "'2.576.55"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("D19").Value = "'4.78"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").Value = "'345.57"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").Value = "'10.41"
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("D22").Value = "'6.09"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").Value = "'59.76"
$ws.Range("E24").Value = "  -0.76%  "
$ws.Range("D25").Value = "'0.419"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("E26").Value = "  -1.07%  "
$ws.Range("D27").Value = "'2.684.07"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("D29").Value = "'0.0₃0842"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").Value = "'7.39"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("D32").Value = "'19.33"
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("D33").Value = "'153.23"
$ws.Range("D34").Value = "'1.56"
$ws.Range("E34").Value = "  -1.39%  "
$ws.Range("D35").Value = "'5.71"
$ws.Range("E35").Value = "  +2.43%  "
$ws.Range("D36").Value = "'3.99"
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("D37").Value = "'1.19"
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("D38").Value = "'0.848"
$ws.Range("E38").Value = "  +7.44%  "
$ws.Range("D39").Value = "'0.850"
$ws.Range("E39").Value = "  -2.34%  "
$ws.Range("D40").Value = "'1.48"
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("D41").Value = "'36.08"
$ws.Range("E41").Value = "  +2.16%  "
$ws.Range("D42").Value = "'3.74"
$ws.Range("D43").Value = "'297.66"
$ws.Range("E43").Value = "  -5.08%  "
$ws.Range("D44").Value = "'0.615"
$ws.Range("E44").Value = "  -3.52%  "
$ws.Range("D45").Value = "'0.0993"
$ws.Range("E45").Value = "  -2.73%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "'0.0557"
$ws.Range("E46").Value = "  -3.11%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").Value = "'0.994"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").Value = "'19.62"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").Value = "'4.82"
$ws.Range("E49").Value = "  -3.51%  "
$ws.Range("D50").Value = "'0.0233"
$ws.Range("E50").Value = "  -2.02%  "
$ws.Range("E51").Value = "  +0.60%  "
